$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new day's row (23/4 -> 2019-04-23) to the daily Page Info log,
# copying the formatting (incl. the date number format) from the row above.
$ws.Range("A62:H62").Copy()
$ws.Range("A63:H63").PasteSpecial(-4122)

$ws.Cells.Item(63, 1).Value = 43578
$ws.Cells.Item(63, 2).Value = 234
$ws.Cells.Item(63, 3).Value = 567
$ws.Cells.Item(63, 4).Value = 1
$ws.Cells.Item(63, 5).Value = 18
$ws.Cells.Item(63, 6).Value = 2476
$ws.Cells.Item(63, 7).Value = 781
$ws.Cells.Item(63, 8).Value = 1987

# Match the saved view's selection state after the edit.
$ws.Range("G66").Select()
